# "Connecting logging to the project"
# The university-profile table is reordered: MATHEMATICS and MEDICINE move
# up (ahead of PHYSICS), while PHYSICS/LINGUISTICS drop down - along with
# each row's average-score / student-count / university-count figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes MATHEMATICS
$ws.Range("B2").Value = "MATHEMATICS"
$ws.Range("C2").Value = 0.0
$ws.Range("D2").Value = 0.0
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = "Казанский Университет Вычислений"

# Row 3 becomes MEDICINE
$ws.Range("B3").Value = "MEDICINE"
$ws.Range("C3").Value = 4.333
$ws.Range("D3").Value = 3.0
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = "Московский Государственный Медицинский Университет, Тамбовский Университет Медицины, Самарский Медицинский Институт"

# Row 4 becomes PHYSICS
$ws.Range("B4").Value = "PHYSICS"
$ws.Range("C4").Value = 4.538
$ws.Range("D4").Value = 8.0
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = "Московский Выдуманный Университет, Московский Придуманный Институт"

# Row 5 becomes LINGUISTICS (counts/score for this profile are unchanged)
$ws.Range("B5").Value = "LINGUISTICS"
$ws.Range("F5").Value = "Воронежский Литературно-Переводческий Университет"

# Refresh the best-fit column widths now that the cell contents changed.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
